# "Changes of New Pre-Prod URL" -- refresh the ShipmentTrackNum / PackageTrackNum
# test values in CheetahProcessing.xlsx (rows 2-22) with a new batch of tracking
# numbers pulled from the new pre-prod environment.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New tracking numbers for column C (ShipmentTrackNum), rows 2-22.
$newValues = @{
    2  = "320018589548"
    3  = "320018589559"
    4  = "320018589581"
    5  = "320018589607"
    6  = "320018589640"
    7  = "320018589662"
    8  = "320018589695"
    9  = "320018589710"
    10 = "320018589743"
    11 = "320018589765"
    12 = "320018589802"
    13 = "320018589824"
    14 = "320018589857"
    15 = "320018589879"
    16 = "320018589905"
    17 = "320018589927"
    18 = "320018589960"
    19 = "320018589982"
    20 = "320018590015"
    21 = "320018590037"
    22 = "320018590060"
}

# Rows where column D (PackageTrackNum) mirrors column C's value.
$dRows = @(5, 6, 7, 13, 14, 15, 16, 17)

# Write every value as a text-typed formula first ("="value"") so the cell
# carries a text result instead of Excel auto-coercing the numeric-looking
# string into a number. Converting that formula in place via Copy +
# PasteSpecial(values) bakes the literal text into the cell without
# touching its existing style.
foreach ($row in $newValues.Keys) {
    $ws.Range("C$row").Formula = '="' + $newValues[$row] + '"'
}
foreach ($row in $dRows) {
    $ws.Range("D$row").Formula = '="' + $newValues[$row] + '"'
}

$cRange = $ws.Range("C2:C22")
$cRange.Copy() | Out-Null
$cRange.PasteSpecial(-4163) | Out-Null

$dRange1 = $ws.Range("D5:D7")
$dRange1.Copy() | Out-Null
$dRange1.PasteSpecial(-4163) | Out-Null

$dRange2 = $ws.Range("D13:D17")
$dRange2.Copy() | Out-Null
$dRange2.PasteSpecial(-4163) | Out-Null

$excel.CutCopyMode = $false
